# Apply the updates described by the commit:
# "Update cfb_weather.xlsx with Timestamp 2024-11-13T10:01:51.609338"
#
# This touches:
#  - the Timestamp column (AK) on the FBS sheet, which holds the same
#    timestamp string in every data row
#  - a handful of recalculated odds / weather values on FBS
#  - a handful of recalculated wind-direction values on FBS and Other

$wb = $excel.ActiveWorkbook

$fbs   = $wb.Worksheets.Item("FBS")
$other = $wb.Worksheets.Item("Other")

# ---------------------------------------------------------------
# Timestamp column: every row (2-51) of column AK on FBS shares the
# same generated timestamp value - refresh them all to the new run.
# ---------------------------------------------------------------
$newTimestamp = "2024-11-13T10:01:51.609338"
$lastRow = $fbs.Cells.Item($fbs.Rows.Count, 1).End(-4162).Row   # xlUp
if ($lastRow -lt 2) { $lastRow = 51 }
$fbs.Range("AK2:AK" + $lastRow).Value = $newTimestamp

# ---------------------------------------------------------------
# FBS sheet numeric / odds updates
# ---------------------------------------------------------------
$fbs.Range("Y2").Value  = 44.5
$fbs.Range("Z2").Value  = -110
$fbs.Range("AE2").Value = -0.04301075268817205

$fbs.Range("Z7").Value  = -115

$fbs.Range("AB32").Value = -7.5
$fbs.Range("AF32").Value = 0

# ---------------------------------------------------------------
# FBS sheet wind direction (wind_dir_fg, column Q) updates
# ---------------------------------------------------------------
$fbs.Range("Q27").Value = "SSE"
$fbs.Range("Q29").Value = "N"
$fbs.Range("Q35").Value = "ENE"
$fbs.Range("Q43").Value = "NNW"
$fbs.Range("Q49").Value = "SSW"
$fbs.Range("Q51").Value = "SW"

# ---------------------------------------------------------------
# Other sheet wind direction (wind_dir_fg, column S) updates
# ---------------------------------------------------------------
$other.Range("S32").Value = "ENE"
$other.Range("S45").Value = "N"
$other.Range("S53").Value = "SSE"
